$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    9  = 5.460199999999998
    13 = 6.659899999999999
    16 = 6.329600000000002
    18 = 6.994499999999992
    20 = 8.7942
    26 = 5.511900000000006
    27 = 6.238900000000001
    29 = 5.174699999999997
    35 = 8.755099999999999
    36 = 9.058500000000006
    45 = 5.038800000000004
    55 = 6.90229999999999
    57 = 5.013299999999995
    69 = 5.307199999999997
    76 = 5.1415
    78 = 9.964500000000001
    82 = 5.455899999999999
    83 = 5.477699999999998
    93 = 5.615999999999998
    97 = 5.488299999999998
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
